# Adds Eetu Pihamäki's 24.10.2018 work-log entry (row 21) to the
# "Eetu Pihamäki" timesheet, and moves the active selection to F22
# (matches the author's post-edit cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eetu Pihamäki")

# New time-log row: date, start time, end time, sprint number, task notes.
# D21 (duration) is a pre-existing shared formula (=C21-B21) that
# recalculates automatically once B21/C21 are populated.
$ws.Range("A21").Value = 43397
$ws.Range("B21").Value = 0.4861111111111111
$ws.Range("C21").Value = 0.8125
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = "30 min ohjauskokousta. 6 h Unix-connectorin asennusta ja konfigurointia. 1h Windowsin 10 liittämistä Windows Server 2016 AD:seen, sekä Markuksen auttamista AD:n kanssa jne. https://github.com/Eetu95/Open-source-IdM-solution/blob/master/Eetun%20muistiinpanoja/Ty%C3%B6t%20-%2024.10.2018.txt"

# The new note text wraps to roughly the same height as other multi-line
# task notes (e.g. row 19), so the row grows to fit it.
$ws.Rows.Item(21).RowHeight = 90

# Author's cursor ended up on F22 after entering the new row.
$ws.Range("F22").Select()
